# Updates market-price-derived columns (currentAveragePrice*, LevePrice*,
# LeveProfit*) across all Leve-profit sheets, refreshing figures pulled from
# the market board data source ("chore: update Sheets via scheduled runner").
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 2082.0833
$ws.Range("I15").Value = 2082.0833
$ws.Range("K15").Value = 6246.249899999999
$ws.Range("M15").Value = -6077.249899999999

# Row 62
$ws.Range("H62").Value = 4315.125
$ws.Range("I62").Value = 3179
$ws.Range("J62").Value = 5451.25
$ws.Range("K62").Value = 3179
$ws.Range("L62").Value = 5451.25
$ws.Range("M62").Value = -2555
$ws.Range("N62").Value = -6699.25

# Row 65
$ws.Range("H65").Value = 4315.125
$ws.Range("I65").Value = 3179
$ws.Range("J65").Value = 5451.25
$ws.Range("K65").Value = 15895
$ws.Range("L65").Value = 27256.25
$ws.Range("M65").Value = -12775
$ws.Range("N65").Value = -33496.25

# Row 106
$ws.Range("H106").Value = 3200
$ws.Range("I106").Value = 1920
$ws.Range("J106").Value = 4000
$ws.Range("K106").Value = 1920
$ws.Range("L106").Value = 4000
$ws.Range("M106").Value = -1289
$ws.Range("N106").Value = -5262

# Row 113
$ws.Range("H113").Value = 4872.8
$ws.Range("I113").Value = 4468.5713
$ws.Range("J113").Value = 5226.5
$ws.Range("K113").Value = 4468.5713
$ws.Range("L113").Value = 5226.5
$ws.Range("M113").Value = -1214.5713
$ws.Range("N113").Value = -11734.5

# Row 116
$ws.Range("H116").Value = 4460.8696
$ws.Range("I116").Value = 3742.2222
$ws.Range("J116").Value = 4922.857
$ws.Range("K116").Value = 3742.2222
$ws.Range("L116").Value = 4922.857
$ws.Range("M116").Value = -300.2222000000002
$ws.Range("N116").Value = -11806.857

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 17978
$ws.Range("I32").Value = 14883.845
$ws.Range("J32").Value = 30635.908
$ws.Range("K32").Value = 14883.845
$ws.Range("L32").Value = 30635.908
$ws.Range("M32").Value = -14596.845
$ws.Range("N32").Value = -31209.908

# Row 37
$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()

# Row 74
$ws.Range("H74").Value = 2064
$ws.Range("I74").Value = 1424.3334
$ws.Range("K74").Value = 1424.3334
$ws.Range("M74").Value = -550.3334

# Row 77
$ws.Range("H77").Value = 2064
$ws.Range("I77").Value = 1424.3334
$ws.Range("K77").Value = 7121.666999999999
$ws.Range("M77").Value = -2753.666999999999

$ws = $wb.Worksheets.Item("BSM")
# Row 59
$ws.Range("H59").Value = 40000
$ws.Range("J59").Value = 40000
$ws.Range("L59").Value = 40000
$ws.Range("N59").Value = -41694

# Row 99
$ws.Range("H99").Value = 3734.2173
$ws.Range("I99").Value = 2856.8235
$ws.Range("J99").Value = 6220.1665
$ws.Range("K99").Value = 2856.8235
$ws.Range("L99").Value = 6220.1665
$ws.Range("M99").Value = -1358.8235
$ws.Range("N99").Value = -9216.166499999999

# Row 105
$ws.Range("H105").Value = 2750.8572
$ws.Range("I105").Value = 2512.2222
$ws.Range("J105").Value = 3180.4
$ws.Range("K105").Value = 2512.2222
$ws.Range("L105").Value = 3180.4
$ws.Range("M105").Value = -765.2222000000002
$ws.Range("N105").Value = -6674.4

$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 16133620
$ws.Range("I58").Value = 2637.4285
$ws.Range("J58").Value = 29417958
$ws.Range("K58").Value = 2637.4285
$ws.Range("L58").Value = 29417958
$ws.Range("M58").Value = -2434.4285
$ws.Range("N58").Value = -29418364

# Row 96
$ws.Range("H96").Value = 22142.285
$ws.Range("J96").Value = 22142.285
$ws.Range("L96").Value = 22142.285
$ws.Range("N96").Value = -27634.285

# Row 122
$ws.Range("H122").Value = 2680
$ws.Range("I122").Value = 2357.375
$ws.Range("J122").Value = 3786.1428
$ws.Range("K122").Value = 7072.125
$ws.Range("L122").Value = 11358.4284
$ws.Range("M122").Value = -4622.125
$ws.Range("N122").Value = -16258.4284

# Row 132
$ws.Range("H132").Value = 2985
$ws.Range("I132").Value = 2902.7693
$ws.Range("J132").Value = 3047.8823
$ws.Range("K132").Value = 8708.3079
$ws.Range("L132").Value = 9143.6469
$ws.Range("M132").Value = -6178.3079
$ws.Range("N132").Value = -14203.6469

# Row 136
$ws.Range("H136").Value = 16133620
$ws.Range("I136").Value = 2637.4285
$ws.Range("J136").Value = 29417958
$ws.Range("K136").Value = 7912.2855
$ws.Range("L136").Value = 88253874
$ws.Range("M136").Value = -5362.2855
$ws.Range("N136").Value = -88258974

$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 1571.4468
$ws.Range("I131").Value = 3124.375
$ws.Range("J131").Value = 1252.8975
$ws.Range("K131").Value = 9373.125
$ws.Range("L131").Value = 3758.6925
$ws.Range("M131").Value = -4333.125
$ws.Range("N131").Value = -13838.6925

# Row 133
$ws.Range("H133").Value = 5443.778
$ws.Range("I133").Value = 4541.4287
$ws.Range("J133").Value = 6018
$ws.Range("K133").Value = 13624.2861
$ws.Range("L133").Value = 18054
$ws.Range("M133").Value = -8564.286100000001
$ws.Range("N133").Value = -28174

$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value = 8098.231
$ws.Range("I122").Value = 769
$ws.Range("K122").Value = 2307
$ws.Range("M122").Value = 143

# Row 123
$ws.Range("H123").Value = 29500
$ws.Range("J123").Value = 29500
$ws.Range("L123").Value = 29500
$ws.Range("N123").Value = -34400

# Row 125
$ws.Range("H125").Value = 21742
$ws.Range("J125").Value = 21742
$ws.Range("L125").Value = 21742
$ws.Range("N125").Value = -26662

# Row 132
$ws.Range("H132").Value = 3418.291
$ws.Range("I132").Value = 3398.2632
$ws.Range("J132").Value = 3463.0588
$ws.Range("K132").Value = 10194.7896
$ws.Range("L132").Value = 10389.1764
$ws.Range("M132").Value = -7664.7896
$ws.Range("N132").Value = -15449.1764

$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 1185.7142
$ws.Range("I46").Value = 903.4483
$ws.Range("J46").Value = 2550
$ws.Range("K46").Value = 903.4483
$ws.Range("L46").Value = 2550
$ws.Range("M46").Value = -715.4483
$ws.Range("N46").Value = -2926

# Row 94
$ws.Range("H94").Value = 25000
$ws.Range("J94").Value = 25000
$ws.Range("L94").Value = 25000
$ws.Range("N94").Value = -26352

# Row 100
$ws.Range("H100").Value = 3921.6365
$ws.Range("I100").Value = 1654.75
$ws.Range("J100").Value = 9966.666999999999
$ws.Range("K100").Value = 1654.75
$ws.Range("L100").Value = 9966.666999999999
$ws.Range("M100").Value = -1113.75
$ws.Range("N100").Value = -11048.667

# Row 122
$ws.Range("H122").Value = 5133.3335
$ws.Range("I122").Value = 4116.6665
$ws.Range("K122").Value = 12349.9995
$ws.Range("M122").Value = -9899.999500000002

# Row 124
$ws.Range("H124").Value = 34000
$ws.Range("J124").Value = 34000
$ws.Range("L124").Value = 34000
$ws.Range("N124").Value = -43820

# Row 132
$ws.Range("H132").Value = 2983.15
$ws.Range("I132").Value = 2242.0908
$ws.Range("K132").Value = 6726.2724
$ws.Range("M132").Value = -4196.2724

# Row 134
$ws.Range("H134").Value = 52400
$ws.Range("J134").Value = 52400
$ws.Range("L134").Value = 52400
$ws.Range("N134").Value = -62540

# Row 141
$ws.Range("H141").Value = 30000
$ws.Range("J141").Value = 30000
$ws.Range("L141").Value = 30000
$ws.Range("N141").Value = -40360

$ws = $wb.Worksheets.Item("WVR")
# Row 22
$ws.Range("H22").Value = 41335.832
$ws.Range("J22").Value = 47603
$ws.Range("L22").Value = 47603
$ws.Range("N22").Value = -48189

# Row 107
$ws.Range("H107").Value = 1908.5834
$ws.Range("I107").Value = 575.6667
$ws.Range("J107").Value = 3241.5
$ws.Range("K107").Value = 1727.0001
$ws.Range("L107").Value = 9724.5
$ws.Range("M107").Value = 192.9999
$ws.Range("N107").Value = -13564.5

# Row 132
$ws.Range("H132").Value = 1972749
$ws.Range("I132").Value = 2443243.5
$ws.Range("J132").Value = 43720.9
$ws.Range("K132").Value = 7329730.5
$ws.Range("L132").Value = 131162.7
$ws.Range("M132").Value = -7327200.5
$ws.Range("N132").Value = -136222.7

# Row 133
$ws.Range("H133").Value = 31857.5
$ws.Range("J133").Value = 31857.5
$ws.Range("L133").Value = 31857.5
$ws.Range("N133").Value = -41977.5
